$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'312.79"
$ws.Range("E2").Value = "'1.33%"
$ws.Range("D3").Value = "'37.58"
$ws.Range("E3").Value = "'0.91%"
$ws.Range("D4").Value = "'5.140"
$ws.Range("E4").Value = "'0.42%"
$ws.Range("E5").Value = "'0.99%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.416"
$ws.Range("E6").Value = "'0.44%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.907"
$ws.Range("E7").Value = "'-3.33%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'8.258"
$ws.Range("E8").Value = "'-0.14%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.850"
$ws.Range("E9").Value = "'-5.96%"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "'0.9210"
$ws.Range("E10").Value = "'-0.42%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.1215"
$ws.Range("E11").Value = "'-7.57%"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1925"
$ws.Range("E12").Value = "'-1.26%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.09133"
$ws.Range("E13").Value = "'1.88%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03283"
$ws.Range("E14").Value = "'-4.55%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09609"
$ws.Range("E15").Value = "'-0.94%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001379"
$ws.Range("E16").Value = "'-0.13%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005714"
$ws.Range("E17").Value = "'-5.67%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.521"
$ws.Range("E18").Value = "'-1.89%"
$ws.Range("E19").Value = "'0.91%"
$ws.Range("D20").Value = "'5.271"
$ws.Range("E20").Value = "'5.15%"
$ws.Range("E21").Value = "'-2.13%"
$ws.Range("D22").Value = "'0.2588"
$ws.Range("E22").Value = "'4.00%"
$ws.Range("E23").Value = "'-0.11%"
$ws.Range("D24").Value = "'0.04363"
$ws.Range("E24").Value = "'0.61%"
$ws.Range("D25").Value = "'0.001248"
$ws.Range("E25").Value = "'2.50%"
$ws.Range("D26").Value = "'0.004313"
$ws.Range("E26").Value = "'-4.51%"
$ws.Range("E27").Value = "'-9.79%"
$ws.Range("D39").Value = "'0.02160"
$ws.Range("E39").Value = "'-5.06%"
$ws.Range("D40").Value = "'0.05118"
$ws.Range("E40").Value = "'1.97%"
$ws.Range("D41").Value = "'0.007461"
$ws.Range("E41").Value = "'-2.51%"
$ws.Range("D42").Value = "'0.1361"
$ws.Range("E42").Value = "'0.63%"
$ws.Range("D43").Value = "'0.008729"
$ws.Range("E43").Value = "'-11.28%"
$ws.Range("D44").Value = "'0.001958"
$ws.Range("E44").Value = "'-2.00%"
$ws.Range("D45").Value = "'0.008634"
$ws.Range("E45").Value = "'2.38%"
$ws.Range("D46").Value = "'0.00006685"
$ws.Range("E46").Value = "'-1.40%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.26%"
$ws.Range("D48").Value = "'0.003350"
$ws.Range("E48").Value = "'11.02%"
$ws.Range("D49").Value = "'0.001199"
$ws.Range("E49").Value = "'-7.78%"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.26%"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.26%"
